$p = $ppt.ActivePresentation
$s1 = $p.Slides.Item(1)

# Remember the style-carrying rectangle so we can copy it onto the
# rebuilt title slide (keeps its <p:style> block intact).
$origRect = $s1.Shapes.Item("Rectangle 1")
$origRect.Copy()

# The target slide is really a "Title Slide" layout (ctrTitle + subTitle
# placeholders) - replace slide 1 with a fresh slide on that layout so the
# placeholders come in with clean shape ids (2, 3) ahead of the pasted
# rectangle (4).
$s1.Delete()
$newSlide = $p.Slides.Add(1, 1)
$s = $p.Slides.Item(1)

# Title placeholder text + color.
$titleRange = $s.Shapes.Item(1).TextFrame.TextRange
$titleRange.Text = "Test Brand Deck"
$titleRange.Font.Color.RGB = 0xC86400

# Subtitle placeholder is left with an empty paragraph (matches target).

# Paste the rectangle that used to live on slide 1, then restyle/move it.
$pasted = $s.Shapes.Paste()
$rect = $s.Shapes.Item($s.Shapes.Count)
$rect.Name = "Rectangle 3"
$rect.Left = 72
$rect.Top = 144
$rect.Width = 288
$rect.Height = 72
$rect.Fill.ForeColor.RGB = 0x3264FF
